$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.531.50"
$ws.Range("E2").Value = "  +1.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.877.87"
$ws.Range("E3").Value = "  +0.79%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7236"
$ws.Range("E5").Value = "  +2.25%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'240.02"
$ws.Range("E6").Value = "  +0.89%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.07865"
$ws.Range("E8").Value = "  -3.94%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.3092"
$ws.Range("E9").Value = "  +1.79%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'25.36"
$ws.Range("E10").Value = "  +8.73%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.08228"
$ws.Range("E11").Value = "  +0.70%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "'1.885.64"
$ws.Range("E12").Value = "  +1.94%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "'0.7272"
$ws.Range("E13").Value = "  +2.66%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.261"
$ws.Range("E14").Value = "  +1.75%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'90.19"
$ws.Range("E15").Value = "  +1.03%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'29.577.68"
$ws.Range("E16").Value = "  +1.16%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "'5.861"
$ws.Range("E17").Value = "  +1.28%  "

# Row 18 - was ShibaInu, now BitcoinCash (rows 18/19 swap)
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'243.00"
$ws.Range("E18").Value = "  +2.68%  "

# Row 19 - was BitcoinCash, now ShibaInu
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007871"
$ws.Range("E19").Value = "  -0.38%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  -0.24%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "'2.135.86"
$ws.Range("E21").Value = "  +2.43%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.12%  "

# Row 23 - BinanceUSD
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'7.786"
$ws.Range("E24").Value = "  +5.21%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1601"
$ws.Range("E25").Value = "  +10.63%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'162.90"
$ws.Range("E26").Value = "  +0.33%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'8.990"
$ws.Range("E27").Value = "  +0.28%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'18.37"
$ws.Range("E28").Value = "  +1.51%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'1.948"
$ws.Range("E29").Value = "  -0.43%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "'1.367"
$ws.Range("E30").Value = "  -4.07%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.487"
$ws.Range("E31").Value = "  +0.19%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.362"
$ws.Range("E32").Value = "  -0.58%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.090"
$ws.Range("E33").Value = "  +0.93%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.05259"
$ws.Range("E34").Value = "  +0.90%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.200"
$ws.Range("E35").Value = "  +2.51%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.7203"
$ws.Range("E36").Value = "  +1.82%  "

# Row 37 - Frax
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38 - HuobiToken
$ws.Range("D38").Value = "'2.673"
$ws.Range("E38").Value = "  -0.03%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.01867"
$ws.Range("E39").Value = "  +0.99%  "

# Row 40 - MXToken
$ws.Range("E40").Value = "  -0.51%  "

# Row 41 - Maker
$ws.Range("D41").Value = "'1.187.67"
$ws.Range("E41").Value = "  +3.84%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "'0.9145"
$ws.Range("E42").Value = "  -1.01%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  +2.32%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "'0.4334"
$ws.Range("E44").Value = "  +1.23%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'72.04"
$ws.Range("E45").Value = "  +2.64%  "

# Row 46 - PaxDollar
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.19%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'103.13"
$ws.Range("E47").Value = "  +0.33%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "'0.5357"
$ws.Range("E48").Value = "  -0.95%  "

# Row 49 - RenderToken
$ws.Range("D49").Value = "'1.780"
$ws.Range("E49").Value = "  +0.26%  "

# Row 50 - was EnergySwap, now SynthetixNetwork (rows 50/51 swap)
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'2.890"
$ws.Range("E50").Value = "  +5.12%  "

# Row 51 - was SynthetixNetwork, now EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.248"
$ws.Range("E51").Value = "  +0.45%  "
